# F_Taxi_UserCaseList_v1.1.xlsx - update
# Signed-off-by: huuphamlc <huuphamlc@gmail.com>
#
# - Move the "Interact with Trip" group label (column C) up by two rows in
#   both affected blocks: C16:C20 -> C14:C20, and C33:C38 -> C31:C38.
# - Swap the D13/D14 "User Case" text (Find Lost Asset / Send Application
#   Feedback were on the wrong rows).
# - Move the active selection to D19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- D13 / D14: swap the two user-case labels -----------------------------
$ws.Range("D13").Value = "Send Application Feedback"
$ws.Range("D14").Value = "Find Lost Asset"

# --- Re-merge the "Interact with Trip" group labels in column C -----------
# Each block grows by two rows at the top (C16:C20 -> C14:C20 and
# C33:C38 -> C31:C38). Do the merge/unmerge first - Excel's Merge() touches
# the alignment of every cell it spans, so we repaint each cell's border
# "slot" (top/middle/bottom of the merged block) afterwards from an
# untouched donor cell with the same look.
$ws.Range("C16:C20").UnMerge()
$ws.Range("C14:C20").Merge()

$ws.Range("C33:C38").UnMerge()
$ws.Range("C31:C38").Merge()

# Donor cells (unaffected by this edit) for each "slot" look:
#   C6  -> top-of-block    (left/right/top border, centered + wrap)
#   C4  -> middle-of-block (left/right border only, centered + wrap)
#   C5  -> bottom-of-block (left/right/bottom border, centered + wrap)

# --- First block: C14:C20 --------------------------------------------------
$ws.Range("C6").Copy()
$ws.Range("C14").PasteSpecial($xlPasteFormats)
$ws.Range("C14").Value = "Interact with Trip"

$ws.Range("C4").Copy()
$ws.Range("C15:C19").PasteSpecial($xlPasteFormats)
$ws.Range("C15:C19").Value = ""

$ws.Range("C5").Copy()
$ws.Range("C20").PasteSpecial($xlPasteFormats)

# --- Second block: C31:C38 -------------------------------------------------
$ws.Range("C6").Copy()
$ws.Range("C31").PasteSpecial($xlPasteFormats)
$ws.Range("C31").Value = "Interact with Trip"

$ws.Range("C4").Copy()
$ws.Range("C32:C37").PasteSpecial($xlPasteFormats)
$ws.Range("C32:C37").Value = ""

$ws.Range("C5").Copy()
$ws.Range("C38").PasteSpecial($xlPasteFormats)

# --- Selection / view ------------------------------------------------------
$ws.Range("D19").Select() | Out-Null

$excel.CutCopyMode = $false
